$d = $word.ActiveDocument

# Remove the "September" milestone sentence (moved to Godot 4 -> dropped the
# "Formal Game Design Document" / greenlight milestone that used to kick off
# the planning in September). The sentence is followed by a manual line
# break (vertical tab) before "In Oktober ..." which also needs removing so
# "In Oktober ..." becomes the first line of the paragraph.
$find = "In September wil ik werken aan een Formal Game Design Document en die gegreenlight hebben door een docent." + [char]11
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
